$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-10-26 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-27 Sunday", 2) | Out-Null
$d.Content.Find.Execute("42+25=", $true, $false, $false, $false, $false, $true, 1, $false, "63-41=", 2) | Out-Null
$d.Content.Find.Execute("28-19=", $true, $false, $false, $false, $false, $true, 1, $false, "56+4=", 2) | Out-Null
$d.Content.Find.Execute("33+64=", $true, $false, $false, $false, $false, $true, 1, $false, "57-18=", 2) | Out-Null
$d.Content.Find.Execute("87-79=", $true, $false, $false, $false, $false, $true, 1, $false, "55+29=", 2) | Out-Null
$d.Content.Find.Execute("23+68=", $true, $false, $false, $false, $false, $true, 1, $false, "38+9=", 2) | Out-Null
$d.Content.Find.Execute("43+48=", $true, $false, $false, $false, $false, $true, 1, $false, "67-26=", 2) | Out-Null
$d.Content.Find.Execute("15+73=", $true, $false, $false, $false, $false, $true, 1, $false, "58-26=", 2) | Out-Null
$d.Content.Find.Execute("85-70=", $true, $false, $false, $false, $false, $true, 1, $false, "26-20=", 2) | Out-Null
$d.Content.Find.Execute("99-30=", $true, $false, $false, $false, $false, $true, 1, $false, "72-36=", 2) | Out-Null
$d.Content.Find.Execute("36+55=", $true, $false, $false, $false, $false, $true, 1, $false, "96-68=", 2) | Out-Null
$d.Content.Find.Execute("41+6=", $true, $false, $false, $false, $false, $true, 1, $false, "6+30=", 2) | Out-Null
$d.Content.Find.Execute("91-52=", $true, $false, $false, $false, $false, $true, 1, $false, "56-49=", 2) | Out-Null
$d.Content.Find.Execute("33-5=", $true, $false, $false, $false, $false, $true, 1, $false, "18+30=", 2) | Out-Null
$d.Content.Find.Execute("78-64=", $true, $false, $false, $false, $false, $true, 1, $false, "58-32=", 2) | Out-Null
$d.Content.Find.Execute("52+25=", $true, $false, $false, $false, $false, $true, 1, $false, "50-29=", 2) | Out-Null
$d.Content.Find.Execute("47-4=", $true, $false, $false, $false, $false, $true, 1, $false, "23+70=", 2) | Out-Null
$d.Content.Find.Execute("16-13=", $true, $false, $false, $false, $false, $true, 1, $false, "36-26=", 2) | Out-Null
$d.Content.Find.Execute("25+59=", $true, $false, $false, $false, $false, $true, 1, $false, "18+1=", 2) | Out-Null
$d.Content.Find.Execute("94-92=", $true, $false, $false, $false, $false, $true, 1, $false, "97-12=", 2) | Out-Null
$d.Content.Find.Execute("98-73=", $true, $false, $false, $false, $false, $true, 1, $false, "38+18=", 2) | Out-Null
$d.Content.Find.Execute("98-13=", $true, $false, $false, $false, $false, $true, 1, $false, "46+9=", 2) | Out-Null
$d.Content.Find.Execute("61-5=", $true, $false, $false, $false, $false, $true, 1, $false, "61+36=", 2) | Out-Null
$d.Content.Find.Execute("69+29=", $true, $false, $false, $false, $false, $true, 1, $false, "60-59=", 2) | Out-Null
$d.Content.Find.Execute("4+72=", $true, $false, $false, $false, $false, $true, 1, $false, "27-2=", 2) | Out-Null
$d.Content.Find.Execute("94-9=", $true, $false, $false, $false, $false, $true, 1, $false, "8+34=", 2) | Out-Null
$d.Content.Find.Execute("61-22=", $true, $false, $false, $false, $false, $true, 1, $false, "0+68=", 2) | Out-Null
$d.Content.Find.Execute("60-23=", $true, $false, $false, $false, $false, $true, 1, $false, "91+4=", 2) | Out-Null
$d.Content.Find.Execute("95-13=", $true, $false, $false, $false, $false, $true, 1, $false, "95-88=", 2) | Out-Null
$d.Content.Find.Execute("97-85=", $true, $false, $false, $false, $false, $true, 1, $false, "7+79=", 2) | Out-Null
$d.Content.Find.Execute("84-4=", $true, $false, $false, $false, $false, $true, 1, $false, "78-62=", 2) | Out-Null
$d.Content.Find.Execute("53+23=", $true, $false, $false, $false, $false, $true, 1, $false, "53+30=", 2) | Out-Null
$d.Content.Find.Execute("84-24=", $true, $false, $false, $false, $false, $true, 1, $false, "67+13=", 2) | Out-Null
$d.Content.Find.Execute("52+26=", $true, $false, $false, $false, $false, $true, 1, $false, "52-25=", 2) | Out-Null
$d.Content.Find.Execute("18+78=", $true, $false, $false, $false, $false, $true, 1, $false, "86-9=", 2) | Out-Null
$d.Content.Find.Execute("63+33=", $true, $false, $false, $false, $false, $true, 1, $false, "28-25=", 2) | Out-Null
$d.Content.Find.Execute("81+13=", $true, $false, $false, $false, $false, $true, 1, $false, "0+3=", 2) | Out-Null
$d.Content.Find.Execute("21-18=", $true, $false, $false, $false, $false, $true, 1, $false, "55-23=", 2) | Out-Null
$d.Content.Find.Execute("95-25=", $true, $false, $false, $false, $false, $true, 1, $false, "24-11=", 2) | Out-Null
$d.Content.Find.Execute("6+67=", $true, $false, $false, $false, $false, $true, 1, $false, "13+17=", 2) | Out-Null
$d.Content.Find.Execute("54+19=", $true, $false, $false, $false, $false, $true, 1, $false, "50+10=", 2) | Out-Null
$d.Content.Find.Execute("90-79=", $true, $false, $false, $false, $false, $true, 1, $false, "83-31=", 2) | Out-Null
$d.Content.Find.Execute("99-78=", $true, $false, $false, $false, $false, $true, 1, $false, "23+51=", 2) | Out-Null
$d.Content.Find.Execute("58-51=", $true, $false, $false, $false, $false, $true, 1, $false, "73-42=", 2) | Out-Null
$d.Content.Find.Execute("82-71=", $true, $false, $false, $false, $false, $true, 1, $false, "53+42=", 2) | Out-Null
$d.Content.Find.Execute("28+31=", $true, $false, $false, $false, $false, $true, 1, $false, "47+39=", 2) | Out-Null
$d.Content.Find.Execute("48-36=", $true, $false, $false, $false, $false, $true, 1, $false, "42-26=", 2) | Out-Null
$d.Content.Find.Execute("99-86=", $true, $false, $false, $false, $false, $true, 1, $false, "72-27=", 2) | Out-Null
$d.Content.Find.Execute("1+80=", $true, $false, $false, $false, $false, $true, 1, $false, "98-7=", 2) | Out-Null
$d.Content.Find.Execute("61-1=", $true, $false, $false, $false, $false, $true, 1, $false, "42-18=", 2) | Out-Null
$d.Content.Find.Execute("52-28=", $true, $false, $false, $false, $false, $true, 1, $false, "65-49=", 2) | Out-Null
$d.Content.Find.Execute("52-27=", $true, $false, $false, $false, $false, $true, 1, $false, "97-44=", 2) | Out-Null
$d.Content.Find.Execute("69+28=", $true, $false, $false, $false, $false, $true, 1, $false, "78-40=", 2) | Out-Null
$d.Content.Find.Execute("84-66=", $true, $false, $false, $false, $false, $true, 1, $false, "57+25=", 2) | Out-Null
$d.Content.Find.Execute("34+35=", $true, $false, $false, $false, $false, $true, 1, $false, "47-13=", 2) | Out-Null
$d.Content.Find.Execute("54-4=", $true, $false, $false, $false, $false, $true, 1, $false, "14+7=", 2) | Out-Null
$d.Content.Find.Execute("64+1=", $true, $false, $false, $false, $false, $true, 1, $false, "66-42=", 2) | Out-Null
$d.Content.Find.Execute("38+47=", $true, $false, $false, $false, $false, $true, 1, $false, "26-10=", 2) | Out-Null
$d.Content.Find.Execute("86-23=", $true, $false, $false, $false, $false, $true, 1, $false, "74-74=", 2) | Out-Null
$d.Content.Find.Execute("51+3=", $true, $false, $false, $false, $false, $true, 1, $false, "4+66=", 2) | Out-Null
$d.Content.Find.Execute("81-75=", $true, $false, $false, $false, $false, $true, 1, $false, "35-23=", 2) | Out-Null
$d.Content.Find.Execute("32+56=", $true, $false, $false, $false, $false, $true, 1, $false, "39-19=", 2) | Out-Null
$d.Content.Find.Execute("26+12=", $true, $false, $false, $false, $false, $true, 1, $false, "67+30=", 2) | Out-Null
$d.Content.Find.Execute("8+83=", $true, $false, $false, $false, $false, $true, 1, $false, "46-43=", 2) | Out-Null
$d.Content.Find.Execute("78-56=", $true, $false, $false, $false, $false, $true, 1, $false, "63-18=", 2) | Out-Null
$d.Content.Find.Execute("64-7=", $true, $false, $false, $false, $false, $true, 1, $false, "4+64=", 2) | Out-Null
$d.Content.Find.Execute("93-56=", $true, $false, $false, $false, $false, $true, 1, $false, "29+48=", 2) | Out-Null
$d.Content.Find.Execute("17+28=", $true, $false, $false, $false, $false, $true, 1, $false, "6+79=", 2) | Out-Null
$d.Content.Find.Execute("34+48=", $true, $false, $false, $false, $false, $true, 1, $false, "6+38=", 2) | Out-Null
$d.Content.Find.Execute("53+18=", $true, $false, $false, $false, $false, $true, 1, $false, "27-19=", 2) | Out-Null
$d.Content.Find.Execute("72-44=", $true, $false, $false, $false, $false, $true, 1, $false, "92-16=", 2) | Out-Null
$d.Content.Find.Execute("54+34=", $true, $false, $false, $false, $false, $true, 1, $false, "21+52=", 2) | Out-Null
$d.Content.Find.Execute("53+27=", $true, $false, $false, $false, $false, $true, 1, $false, "99-70=", 2) | Out-Null
$d.Content.Find.Execute("28+17=", $true, $false, $false, $false, $false, $true, 1, $false, "68+22=", 2) | Out-Null
$d.Content.Find.Execute("47-46=", $true, $false, $false, $false, $false, $true, 1, $false, "92-33=", 2) | Out-Null
$d.Content.Find.Execute("30+8=", $true, $false, $false, $false, $false, $true, 1, $false, "72-41=", 2) | Out-Null
$d.Content.Find.Execute("53-16=", $true, $false, $false, $false, $false, $true, 1, $false, "46+52=", 2) | Out-Null
$d.Content.Find.Execute("71+10=", $true, $false, $false, $false, $false, $true, 1, $false, "56+1=", 2) | Out-Null
$d.Content.Find.Execute("13+4=", $true, $false, $false, $false, $false, $true, 1, $false, "85-10=", 2) | Out-Null
$d.Content.Find.Execute("5+55=", $true, $false, $false, $false, $false, $true, 1, $false, "70-25=", 2) | Out-Null
$d.Content.Find.Execute("82-63=", $true, $false, $false, $false, $false, $true, 1, $false, "4+41=", 2) | Out-Null
$d.Content.Find.Execute("40+47=", $true, $false, $false, $false, $false, $true, 1, $false, "71+20=", 2) | Out-Null
$d.Content.Find.Execute("81-58=", $true, $false, $false, $false, $false, $true, 1, $false, "33+53=", 2) | Out-Null
$d.Content.Find.Execute("85-29=", $true, $false, $false, $false, $false, $true, 1, $false, "4+30=", 2) | Out-Null
$d.Content.Find.Execute("1+12=", $true, $false, $false, $false, $false, $true, 1, $false, "24+9=", 2) | Out-Null
$d.Content.Find.Execute("18+36=", $true, $false, $false, $false, $false, $true, 1, $false, "11+12=", 2) | Out-Null
$d.Content.Find.Execute("77+0=", $true, $false, $false, $false, $false, $true, 1, $false, "15+61=", 2) | Out-Null
$d.Content.Find.Execute("48+3=", $true, $false, $false, $false, $false, $true, 1, $false, "81-73=", 2) | Out-Null
$d.Content.Find.Execute("91-89=", $true, $false, $false, $false, $false, $true, 1, $false, "21+43=", 2) | Out-Null
$d.Content.Find.Execute("86-45=", $true, $false, $false, $false, $false, $true, 1, $false, "91-63=", 2) | Out-Null
$d.Content.Find.Execute("29+45=", $true, $false, $false, $false, $false, $true, 1, $false, "38+11=", 2) | Out-Null
$d.Content.Find.Execute("76-23=", $true, $false, $false, $false, $false, $true, 1, $false, "41+23=", 2) | Out-Null
$d.Content.Find.Execute("94-18=", $true, $false, $false, $false, $false, $true, 1, $false, "48-11=", 2) | Out-Null
$d.Content.Find.Execute("29+56=", $true, $false, $false, $false, $false, $true, 1, $false, "51-15=", 2) | Out-Null
$d.Content.Find.Execute("17+12=", $true, $false, $false, $false, $false, $true, 1, $false, "12+13=", 2) | Out-Null
$d.Content.Find.Execute("79-51=", $true, $false, $false, $false, $false, $true, 1, $false, "81-40=", 2) | Out-Null
$d.Content.Find.Execute("96-18=", $true, $false, $false, $false, $false, $true, 1, $false, "87-71=", 2) | Out-Null
$d.Content.Find.Execute("2+46=", $true, $false, $false, $false, $false, $true, 1, $false, "14+4=", 2) | Out-Null
$d.Content.Find.Execute("53+44=", $true, $false, $false, $false, $false, $true, 1, $false, "73-30=", 2) | Out-Null
$d.Content.Find.Execute("79+12=", $true, $false, $false, $false, $false, $true, 1, $false, "18+54=", 2) | Out-Null
$d.Content.Find.Execute("42-27=", $true, $false, $false, $false, $false, $true, 1, $false, "93-42=", 2) | Out-Null
